$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.319.12'
$ws.Range('E2').Value = '  +2.30%  '
$ws.Range('D3').Value = '3.585.81'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.35'
$ws.Range('E5').Value = '  +3.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.77'
$ws.Range('E6').Value = '  +19.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '652.38'
$ws.Range('E7').Value = '  -0.61%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.430'
$ws.Range('E8').Value = '  +7.34%  '
$ws.Range('E10').Value = '  +4.07%  '
$ws.Range('D11').Value = '3.582.71'
$ws.Range('E11').Value = '  +0.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '44.19'
$ws.Range('E12').Value = '  +4.05%  '
$ws.Range('E13').Value = '  +0.97%  '
$ws.Range('E14').Value = '  +0.51%  '
$ws.Range('D15').Value = '4.251.64'
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('D16').Value = '97.081.76'
$ws.Range('E16').Value = '  +2.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000261'
$ws.Range('D18').Value = '3.587.49'
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.77'
$ws.Range('E19').Value = '  -1.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.60'
$ws.Range('E20').Value = '  -0.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.11'
$ws.Range('E21').Value = '  +1.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.537'
$ws.Range('E22').Value = '  +11.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '518.81'
$ws.Range('E23').Value = '  +1.97%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('E25').Value = '  +4.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.94'
$ws.Range('E26').Value = '  -0.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '103.81'
$ws.Range('E27').Value = '  +13.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '13.17'
$ws.Range('E28').Value = '  +2.50%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.176'
$ws.Range('E29').Value = '  +21.55%  '
$ws.Range('B30').Value = 'WrappedeETH'
$ws.Range('C30').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D30').Value = '3.778.50'
$ws.Range('E30').Value = '  +0.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.99'
$ws.Range('E31').Value = '  -1.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.01'
$ws.Range('E32').Value = '  +3.79%  '
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  +5.75%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.92'
$ws.Range('E36').Value = '  -0.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.578'
$ws.Range('E37').Value = '  +2.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '616.24'
$ws.Range('E38').Value = '  +2.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.68'
$ws.Range('E39').Value = '  +1.14%  '
$ws.Range('E40').Value = '  -5.35%  '
$ws.Range('E41').Value = '  +1.50%  '
$ws.Range('E42').Value = '  +4.76%  '
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.931'
$ws.Range('E44').Value = '  +2.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.10'
$ws.Range('E45').Value = '  +6.20%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.439'
$ws.Range('E46').Value = '  +40.35%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0445'
$ws.Range('E47').Value = '  +6.95%  '
$ws.Range('E48').Value = '  +0.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.64'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.65'
$ws.Range('E50').Value = '  +5.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.30'
$ws.Range('E51').Value = '  +7.53%  '
